# Weekly update: insert three new report rows (new reporting date 2021-10-05)
# at the top of the Brócoli / Lo Valledor data block (row 377), pushing the
# existing rows down by three (old 377-437 -> new 380-440).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 377 (existing rows 377:437 shift down to 380:440)
$ws.Rows("377:379").Insert()

# ---- New row 377 ----
$ws.Cells.Item(377, 1).Value = 6
$ws.Cells.Item(377, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(377, 3).Value = "Metropolitana"
$ws.Cells.Item(377, 4).Value = 44474
$ws.Cells.Item(377, 5).Value = 13
$ws.Cells.Item(377, 6).Value = 100112023
$ws.Cells.Item(377, 7).Value = "Brócoli"
$ws.Cells.Item(377, 8).Value = "Sin especificar"
$ws.Cells.Item(377, 9).Value = "Primera"
$ws.Cells.Item(377, 10).Value = 15500
$ws.Cells.Item(377, 11).Value = 450
$ws.Cells.Item(377, 12).Value = 600
$ws.Cells.Item(377, 13).Value = 524
$ws.Cells.Item(377, 14).Value = "`$/unidad"
$ws.Cells.Item(377, 15).Value = "Región Metropolitana"
$ws.Cells.Item(377, 16).Value = 524
$ws.Cells.Item(377, 17).Value = 1
$ws.Cells.Item(377, 18).Value = "Hortaliza"

# ---- New row 378 ----
$ws.Cells.Item(378, 1).Value = 6
$ws.Cells.Item(378, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(378, 3).Value = "Metropolitana"
$ws.Cells.Item(378, 4).Value = 44474
$ws.Cells.Item(378, 5).Value = 13
$ws.Cells.Item(378, 6).Value = 100112023
$ws.Cells.Item(378, 7).Value = "Brócoli"
$ws.Cells.Item(378, 8).Value = "Sin especificar"
$ws.Cells.Item(378, 9).Value = "Primera"
$ws.Cells.Item(378, 10).Value = 6700
$ws.Cells.Item(378, 11).Value = 600
$ws.Cells.Item(378, 12).Value = 700
$ws.Cells.Item(378, 13).Value = 651
$ws.Cells.Item(378, 14).Value = "`$/unidad"
$ws.Cells.Item(378, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(378, 16).Value = 651
$ws.Cells.Item(378, 17).Value = 1
$ws.Cells.Item(378, 18).Value = "Hortaliza"

# ---- New row 379 ----
$ws.Cells.Item(379, 1).Value = 6
$ws.Cells.Item(379, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(379, 3).Value = "Metropolitana"
$ws.Cells.Item(379, 4).Value = 44474
$ws.Cells.Item(379, 5).Value = 13
$ws.Cells.Item(379, 6).Value = 100112023
$ws.Cells.Item(379, 7).Value = "Brócoli"
$ws.Cells.Item(379, 8).Value = "Sin especificar"
$ws.Cells.Item(379, 9).Value = "Segunda"
$ws.Cells.Item(379, 10).Value = 5500
$ws.Cells.Item(379, 11).Value = 300
$ws.Cells.Item(379, 12).Value = 350
$ws.Cells.Item(379, 13).Value = 325
$ws.Cells.Item(379, 14).Value = "`$/unidad"
$ws.Cells.Item(379, 15).Value = "Región Metropolitana"
$ws.Cells.Item(379, 16).Value = 325
$ws.Cells.Item(379, 17).Value = 1
$ws.Cells.Item(379, 18).Value = "Hortaliza"
